$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3480.825
$ws.Range("I80").Value = 788.86664
$ws.Range("J80").Value = 5096
$ws.Range("K80").Value = 2366.59992
$ws.Range("L80").Value = 15288
$ws.Range("M80").Value = -1368.59992
$ws.Range("N80").Value = -17284

$ws.Range("H83").Value = 3480.825
$ws.Range("I83").Value = 788.86664
$ws.Range("J83").Value = 5096
$ws.Range("K83").Value = 7099.79976
$ws.Range("L83").Value = 45864
$ws.Range("M83").Value = -2107.79976
$ws.Range("N83").Value = -55848

$ws.Range("H98").Value = 1564.9286
$ws.Range("I98").Value = 1699.3889
$ws.Range("J98").Value = 758.1667
$ws.Range("K98").Value = 1699.3889
$ws.Range("L98").Value = 758.1667
$ws.Range("M98").Value = -201.3888999999999
$ws.Range("N98").Value = -3754.1667

$ws.Range("H122").Value = 1564.9286
$ws.Range("I122").Value = 1699.3889
$ws.Range("J122").Value = 758.1667
$ws.Range("K122").Value = 5098.1667
$ws.Range("L122").Value = 2274.5001
$ws.Range("M122").Value = -2648.1667
$ws.Range("N122").Value = -7174.5001

$ws.Range("H132").Value = 1339.3383
$ws.Range("I132").Value = 1156.9672
$ws.Range("K132").Value = 3470.9016
$ws.Range("M132").Value = -940.9016000000001

$ws.Range("H138").Value = 3959.361
$ws.Range("I138").Value = 2866
$ws.Range("J138").Value = 4271.75
$ws.Range("K138").Value = 8598
$ws.Range("L138").Value = 12815.25
$ws.Range("M138").Value = -3458
$ws.Range("N138").Value = -23095.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 6692
$ws.Range("I88").Value = 14279.5
$ws.Range("J88").Value = 2356.2856
$ws.Range("K88").Value = 14279.5
$ws.Range("L88").Value = 2356.2856
$ws.Range("M88").Value = -13873.5
$ws.Range("N88").Value = -3168.2856

$ws.Range("H91").Value = 6692
$ws.Range("I91").Value = 14279.5
$ws.Range("J91").Value = 2356.2856
$ws.Range("K91").Value = 14279.5
$ws.Range("L91").Value = 2356.2856
$ws.Range("M91").Value = -12875.5
$ws.Range("N91").Value = -5164.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2370.7856
$ws.Range("I107").Value = 1987.8889
$ws.Range("K107").Value = 1987.8889
$ws.Range("M107").Value = -67.88889999999992

$ws.Range("H134").Value = 27451.074
$ws.Range("I134").Value = 2428.6
$ws.Range("K134").Value = 7285.799999999999
$ws.Range("M134").Value = -4750.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2398.8525
$ws.Range("I31").Value = 1728.32
$ws.Range("J31").Value = 5446.727
$ws.Range("K31").Value = 1728.32
$ws.Range("L31").Value = 5446.727
$ws.Range("M31").Value = -1433.32
$ws.Range("N31").Value = -6036.727

$ws.Range("H34").Value = 2398.8525
$ws.Range("I34").Value = 1728.32
$ws.Range("J34").Value = 5446.727
$ws.Range("K34").Value = 1728.32
$ws.Range("L34").Value = 5446.727
$ws.Range("M34").Value = -1526.32
$ws.Range("N34").Value = -5850.727

$ws.Range("H105").Value = 346.27274
$ws.Range("I105").Value = 359.8
$ws.Range("K105").Value = 359.8
$ws.Range("M105").Value = 1387.2

$ws.Range("H122").Value = 9845.5
$ws.Range("I122").Value = 5911.1113
$ws.Range("K122").Value = 17733.3339
$ws.Range("M122").Value = -15283.3339

$ws.Range("H132").Value = 2192.1296
$ws.Range("I132").Value = 1968.0435
$ws.Range("J132").Value = 2358.3872
$ws.Range("K132").Value = 5904.1305
$ws.Range("L132").Value = 7075.1616
$ws.Range("M132").Value = -3374.1305
$ws.Range("N132").Value = -12135.1616

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 15385348
$ws.Range("I18").Value = 22222470
$ws.Range("J18").Value = 1825
$ws.Range("K18").Value = 66667410
$ws.Range("L18").Value = 5475
$ws.Range("M18").Value = -66667241
$ws.Range("N18").Value = -5813

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H110").Value = 4997.0586
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 4997.0586
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 14991.1758
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -23171.1758

$ws.Range("H138").Value = 6709.909
$ws.Range("I138").Value = 11156
$ws.Range("J138").Value = 3004.8333
$ws.Range("K138").Value = 33468
$ws.Range("L138").Value = 9014.499899999999
$ws.Range("M138").Value = -28328
$ws.Range("N138").Value = -19294.4999

$ws.Range("H139").Value = 2074106.8
$ws.Range("I139").Value = 3914973.2
$ws.Range("J139").Value = 3132.0625
$ws.Range("K139").Value = 11744919.6
$ws.Range("L139").Value = 9396.1875
$ws.Range("M139").Value = -11739779.6
$ws.Range("N139").Value = -19676.1875

$ws.Range("H141").Value = 3837.1667
$ws.Range("I141").Value = 2922.9
$ws.Range("K141").Value = 8768.700000000001
$ws.Range("M141").Value = -3588.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3210.7917
$ws.Range("I102").Value = 2881.3428
$ws.Range("K102").Value = 2881.3428
$ws.Range("M102").Value = -1259.3428

$ws.Range("H122").Value = 8952.25
$ws.Range("I122").Value = 26000
$ws.Range("K122").Value = 78000
$ws.Range("M122").Value = -75550

$ws.Range("H126").Value = 3083.3635
$ws.Range("I126").Value = 2025
$ws.Range("J126").Value = 3688.1428
$ws.Range("K126").Value = 6075
$ws.Range("L126").Value = 11064.4284
$ws.Range("M126").Value = -3605
$ws.Range("N126").Value = -16004.4284

$ws.Range("H132").Value = 35059.766
$ws.Range("I132").Value = 64257
$ws.Range("J132").Value = 9106.666999999999
$ws.Range("K132").Value = 192771
$ws.Range("L132").Value = 27320.001
$ws.Range("M132").Value = -190241
$ws.Range("N132").Value = -32380.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3280.3572
$ws.Range("I7").Value = 2998.7222
$ws.Range("J7").Value = 3787.3
$ws.Range("K7").Value = 2998.7222
$ws.Range("L7").Value = 3787.3
$ws.Range("M7").Value = -2886.7222
$ws.Range("N7").Value = -4011.3

$ws.Range("H40").Value = 4590.6665
$ws.Range("I40").Value = 4337.75
$ws.Range("J40").Value = 5400
$ws.Range("K40").Value = 4337.75
$ws.Range("L40").Value = 5400
$ws.Range("M40").Value = -4201.75
$ws.Range("N40").Value = -5672

$ws.Range("H82").Value = 2542.3572
$ws.Range("I82").Value = 1470
$ws.Range("J82").Value = 3614.7144
$ws.Range("K82").Value = 1470
$ws.Range("L82").Value = 3614.7144
$ws.Range("M82").Value = -1109
$ws.Range("N82").Value = -4336.7144

$ws.Range("H85").Value = 2542.3572
$ws.Range("I85").Value = 1470
$ws.Range("J85").Value = 3614.7144
$ws.Range("K85").Value = 1470
$ws.Range("L85").Value = 3614.7144
$ws.Range("M85").Value = -222
$ws.Range("N85").Value = -6110.7144

$ws.Range("H126").Value = 3280.3572
$ws.Range("I126").Value = 2998.7222
$ws.Range("J126").Value = 3787.3
$ws.Range("K126").Value = 8996.1666
$ws.Range("L126").Value = 11361.9
$ws.Range("M126").Value = -6526.1666
$ws.Range("N126").Value = -16301.9

$ws.Range("H132").Value = 3240.9395
$ws.Range("I132").Value = 2967.125
$ws.Range("K132").Value = 8901.375
$ws.Range("M132").Value = -6371.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3873.2903
$ws.Range("I122").Value = 2304.0715
$ws.Range("J122").Value = 5165.5884
$ws.Range("K122").Value = 6912.2145
$ws.Range("L122").Value = 15496.7652
$ws.Range("M122").Value = -4462.2145
$ws.Range("N122").Value = -20396.7652

$ws.Range("H132").Value = 2872.4524
$ws.Range("J132").Value = 5641.5557
$ws.Range("L132").Value = 16924.6671
$ws.Range("M132").Value = -2872.7272
$ws.Range("N132").Value = -21984.6671
